# Add a new "Partition coefficients" worksheet between "Male " and "Female",
# populate it with the partition-coefficient table for the three compounds
# (Cinnamaldehyde, Cinnamyl Alcohol, Benzaldehyde) and turn the range into a
# native Excel Table ("Tabel1"), matching the committed workbook.

$wb = $excel.ActiveWorkbook

$maleSheet = $wb.Worksheets.Item("Male ")

# Insert the new sheet right after "Male " (i.e. before "Female").
$ws = $wb.Worksheets.Add($null, $maleSheet)
$ws.Name = "Partition coefficients"

# Column widths (matches the authored sheet).
$ws.Columns.Item(1).ColumnWidth = 30.21875
$ws.Columns.Item(2).ColumnWidth = 11.88671875
$ws.Columns.Item(3).ColumnWidth = 22.77734375

# Header row.
$ws.Range("A1").Value = "Partition coefficients "
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Variable name in R"

# Section headers (bold).
$ws.Range("A2").Value = "Cinnamaldehyde "

$ws.Range("A10").Value = "Cinnamyl Alcohol"
$ws.Range("A10").Font.Bold = $true

$ws.Range("A17").Value = "Benzaldehyde "
$ws.Range("A17").Font.Bold = $true

# Data rows: row, Label (A), Value (B), R variable name (C)
$rowsData = @(
    @(3, "Fat:Blood ", 1.62, "P_F"),
    @(4, "liver:Blood ", 0.59, "P_L"),
    @(5, " Small intestine:Blood ", 0.59, "P_SI"),
    @(6, "Richly perfused tissues:Blood", 0.59, "P_RP"),
    @(7, "Slowly perfused tissues:Blood ", 0.78, "P_SP"),
    @(8, "Blood:Air ", "1,25*10^5", "P_PB"),
    @(9, "lung:Blood", 0.59, "P_Pu"),
    @(11, "Fat:Blood ", 1.64, "P_OH_F "),
    @(12, "liver:Blood ", 0.59, "P_OH_L "),
    @(13, " Small intestine:Blood ", 0.59, "P_OH_SI"),
    @(14, "Richly perfused tissues:Blood", 0.59, "P_OH_RP"),
    @(15, "Slowly perfused tissues:Blood ", 0.78, "P_OH_SP"),
    @(16, "lung:Blood", 0.59, "P_OH_Pu"),
    @(18, "Fat:Blood ", 1.51, "P_F_Benz"),
    @(19, "liver:Blood ", 0.59, "P_L_Benz"),
    @(20, " Small intestine:Blood ", 0.59, "P_SI_Benz"),
    @(21, "Richly perfused tissues:Blood", 0.59, "P_RP_Benz"),
    @(22, "Slowly perfused tissues:Blood ", 0.78, "P_SP_Benz"),
    @(23, "Blood:Air ", "1,66*10^4", "P_PB_Benz"),
    @(24, "lung:Blood", 0.59, "P_Pu")
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $label = $row[1]
    $value = $row[2]
    $varName = $row[3]

    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $value
    $ws.Cells.Item($r, 3).Value = $varName
}

# "Blood:Air" value cells are right-aligned text in the source workbook.
$ws.Range("B8").HorizontalAlignment = -4152
$ws.Range("B23").HorizontalAlignment = -4152

# Selection the author left on this sheet.
$ws.Range("F16").Select()

# Turn the populated range into a native Excel table.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:C24"), $null, 1)
$lo.Name = "Tabel1"
$lo.TableStyle = "TableStyleMedium7"

# Make the new sheet the active tab, matching the saved workbook view.
$ws.Activate()
